$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E2").Value = 118
$ws.Range("F2").Value = 87
$ws.Range("H2").Value = 93
$ws.Range("E3").Value = 47
$ws.Range("F3").Value = 36
$ws.Range("H3").Value = 37
$ws.Range("F4").Value = 33
$ws.Range("H4").Value = 45
$ws.Range("E7").Value = 50
$ws.Range("F7").Value = 34
$ws.Range("H7").Value = 38
$ws.Range("E10").Value = 767
$ws.Range("F10").Value = 457
$ws.Range("H10").Value = 552
$ws.Range("E11").Value = 508
$ws.Range("F11").Value = 314
$ws.Range("H11").Value = 379
$ws.Range("E12").Value = 793
$ws.Range("F12").Value = 506
$ws.Range("H12").Value = 592
$ws.Range("E13").Value = 179
$ws.Range("F13").Value = 106
$ws.Range("H13").Value = 140
$ws.Range("E14").Value = 154
$ws.Range("E15").Value = 220
$ws.Range("F15").Value = 105
$ws.Range("H15").Value = 156
$ws.Range("E16").Value = 245
$ws.Range("F16").Value = 146
$ws.Range("H16").Value = 194
$ws.Range("E17").Value = 135
$ws.Range("F17").Value = 75
$ws.Range("H17").Value = 99
$ws.Range("E18").Value = 66
$ws.Range("F18").Value = 39
$ws.Range("H18").Value = 56
$ws.Range("E20").Value = 106
$ws.Range("F20").Value = 52
$ws.Range("H20").Value = 89
$ws.Range("E21").Value = 154
$ws.Range("E22").Value = 208
$ws.Range("F22").Value = 123
$ws.Range("H22").Value = 165
$ws.Range("E23").Value = 243
$ws.Range("F23").Value = 128
$ws.Range("H23").Value = 180
$ws.Range("E24").Value = 301
$ws.Range("F24").Value = 177
$ws.Range("H24").Value = 207
$ws.Range("E25").Value = 357
$ws.Range("F25").Value = 208
$ws.Range("H25").Value = 268
$ws.Range("E26").Value = 234
$ws.Range("F26").Value = 147
$ws.Range("H26").Value = 172
$ws.Range("E27").Value = 420
$ws.Range("F27").Value = 246
$ws.Range("H27").Value = 328
$ws.Range("E28").Value = 242
$ws.Range("F28").Value = 121
$ws.Range("H28").Value = 173
$ws.Range("E29").Value = 204
$ws.Range("F29").Value = 127
$ws.Range("H29").Value = 168
$ws.Range("E30").Value = 277
$ws.Range("F30").Value = 179
$ws.Range("H30").Value = 232
$ws.Range("F31").Value = 41
$ws.Range("H31").Value = 68
$ws.Range("E32").Value = 227
$ws.Range("F32").Value = 154
$ws.Range("H32").Value = 192
$ws.Range("E33").Value = 364
$ws.Range("F33").Value = 197
$ws.Range("H33").Value = 288
$ws.Range("E34").Value = 274
$ws.Range("F34").Value = 196
$ws.Range("H34").Value = 234
$ws.Range("E35").Value = 193
$ws.Range("F35").Value = 139
$ws.Range("H35").Value = 166
$ws.Range("E36").Value = 93
$ws.Range("F36").Value = 60
$ws.Range("H36").Value = 70
$ws.Range("E37").Value = 207
$ws.Range("F37").Value = 119
$ws.Range("H37").Value = 155
$ws.Range("E38").Value = 113
$ws.Range("F38").Value = 74
$ws.Range("H38").Value = 91
$ws.Range("E39").Value = 213
$ws.Range("F39").Value = 113
$ws.Range("H39").Value = 164
$ws.Range("E40").Value = 325
$ws.Range("F40").Value = 178
$ws.Range("H40").Value = 258
$ws.Range("E41").Value = 466
$ws.Range("F41").Value = 250
$ws.Range("H41").Value = 342
$ws.Range("E42").Value = 498
$ws.Range("F42").Value = 311
$ws.Range("H42").Value = 372
$ws.Range("E43").Value = 155
$ws.Range("F43").Value = 94
$ws.Range("H43").Value = 121
$ws.Range("E44").Value = 408
$ws.Range("F44").Value = 232
$ws.Range("H44").Value = 300
$ws.Range("F45").Value = 120
$ws.Range("H45").Value = 159
$ws.Range("E46").Value = 415
$ws.Range("F46").Value = 251
$ws.Range("H46").Value = 315
$ws.Range("E47").Value = 587
$ws.Range("F47").Value = 345
$ws.Range("H47").Value = 437
$ws.Range("E48").Value = 296
$ws.Range("F48").Value = 151
$ws.Range("H48").Value = 195
$ws.Range("E49").Value = 356
$ws.Range("F49").Value = 185
$ws.Range("H49").Value = 272
$ws.Range("E50").Value = 299
$ws.Range("F50").Value = 174
$ws.Range("H50").Value = 247
$ws.Range("E51").Value = 270
$ws.Range("F51").Value = 143
$ws.Range("H51").Value = 217
